$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column E (Foreign Born (%)) values and one F value (Missing For. Born count) per diff
$ws.Range("E2").Value = 0.823268921095008
$ws.Range("E3").Value = 4.51781012841525
$ws.Range("E4").Value = 16.4758181108709
$ws.Range("E5").Value = 24.85192000257
$ws.Range("E6").Value = 0.346389942857143
$ws.Range("E7").Value = 5.5916040259553
$ws.Range("F7").Value = 6
$ws.Range("E8").Value = 23.1284109339484
$ws.Range("E9").Value = 47.905288968098
$ws.Range("E14").Value = 0.567888133728462
$ws.Range("E15").Value = 1.20445965009235
$ws.Range("E16").Value = 4.69300321111288
$ws.Range("E17").Value = 36.4188162835123
